$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$textCells = @("D5","D6","D7","D8","D9","D10","D14","D15","D17","D18","D19","D22","D23","D25","D26","D27","D30","D33","D34","D36","D38","D41","D42","D43","D47","D48","D49","D51")
foreach ($addr in $textCells) { $ws.Range($addr).NumberFormat = "@" }

$ws.Range("D2").Value = "34.038.52"
$ws.Range("E2").Value = "  -1.38%  "
$ws.Range("D3").Value = "1.784.47"
$ws.Range("E3").Value = "  -3.04%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "224.06"
$ws.Range("E5").Value = "  -1.15%  "
$ws.Range("D6").Value = "0.550"
$ws.Range("E6").Value = "  -0.96%  "
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").Value = "32.35"
$ws.Range("E8").Value = "  -0.31%  "
$ws.Range("D9").Value = "0.284"
$ws.Range("E9").Value = "  -3.98%  "
$ws.Range("D10").Value = "0.0704"
$ws.Range("E10").Value = "  -2.35%  "
$ws.Range("E11").Value = "  -0.30%  "
$ws.Range("D12").Value = "2.041.13"
$ws.Range("E12").Value = "  -3.04%  "
$ws.Range("D13").Value = "1.779.64"
$ws.Range("E13").Value = "  -3.39%  "
$ws.Range("D14").Value = "10.81"
$ws.Range("E14").Value = "  -1.70%  "
$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").Value = "0.622"
$ws.Range("E15").Value = "  -4.49%  "
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "34.029.85"
$ws.Range("E16").Value = "  -1.53%  "
$ws.Range("D17").Value = "4.15"
$ws.Range("E17").Value = "  -5.13%  "
$ws.Range("D18").Value = "67.75"
$ws.Range("D19").Value = "243.28"
$ws.Range("E19").Value = "  -3.93%  "
$ws.Range("D20").Value = "0.0₃0782"
$ws.Range("E20").Value = "  -3.25%  "
$ws.Range("E21").Value = "  +0.18%  "
$ws.Range("D22").Value = "10.73"
$ws.Range("E22").Value = "  -5.08%  "
$ws.Range("D23").Value = "4.09"
$ws.Range("E23").Value = "  -5.17%  "
$ws.Range("E24").Value = "  -2.47%  "
$ws.Range("D25").Value = "159.67"
$ws.Range("E25").Value = "  -1.38%  "
$ws.Range("D26").Value = "16.29"
$ws.Range("E26").Value = "  -3.27%  "
$ws.Range("D27").Value = "7.03"
$ws.Range("E27").Value = "  -3.34%  "
$ws.Range("E28").Value = "  -2.74%  "
$ws.Range("E29").Value = "  +0.09%  "
$ws.Range("D30").Value = "0.0513"
$ws.Range("E30").Value = "  -4.75%  "
$ws.Range("E31").Value = "  -0.07%  "
$ws.Range("E32").Value = "  -4.44%  "
$ws.Range("D33").Value = "3.49"
$ws.Range("E33").Value = "  -4.13%  "
$ws.Range("D34").Value = "1.82"
$ws.Range("E34").Value = "  -7.35%  "
$ws.Range("D35").Value = "1.392.41"
$ws.Range("E35").Value = "  -4.87%  "
$ws.Range("D36").Value = "0.644"
$ws.Range("E36").Value = "  -1.86%  "
$ws.Range("E37").Value = "  -3.37%  "
$ws.Range("D38").Value = "0.0186"
$ws.Range("E38").Value = "  -4.49%  "
$ws.Range("E39").Value = "  -0.89%  "
$ws.Range("E40").Value = "  +1.81%  "
$ws.Range("D41").Value = "2.70"
$ws.Range("E41").Value = "  -3.49%  "
$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D42").Value = "0.911"
$ws.Range("E42").Value = "  -7.39%  "
$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").Value = "78.18"
$ws.Range("E43").Value = "  -5.89%  "
$ws.Range("E44").Value = "  +12.70%  "
$ws.Range("E45").Value = "  +1.03%  "
$ws.Range("E46").Value = "  -0.21%  "
$ws.Range("D47").Value = "107.01"
$ws.Range("D48").Value = "5.86"
$ws.Range("E48").Value = "  -4.78%  "
$ws.Range("D49").Value = "12.31"
$ws.Range("E49").Value = "  -0.57%  "
$ws.Range("D50").Value = "1.941.24"
$ws.Range("E50").Value = "  -2.93%  "
$ws.Range("D51").Value = "0.999"
